$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @"
2|46059|0|06.02.202615
3|46059|0|06.02.202616
4|46059|0.123|06.02.202617
5|46059|0.063|06.02.202618
6|46059|0|06.02.202619
7|46059|0|06.02.202620
8|46059|0|06.02.202621
9|46059|0|06.02.202622
10|46059|0|06.02.202623
11|46059|0|06.02.202624
12|46060|0|07.02.20261
13|46060|0|07.02.20262
14|46060|0|07.02.20263
15|46060|0|07.02.20264
16|46060|0|07.02.20265
17|46060|0|07.02.20266
18|46060|0|07.02.20267
19|46060|0|07.02.20268
20|46060|0|07.02.20269
21|46060|0.121|07.02.202610
22|46060|0.517|07.02.202611
23|46060|0.529|07.02.202612
24|46060|0.671|07.02.202613
25|46060|0.71|07.02.202614
26|46060|0.71|07.02.202615
27|46060|0.584|07.02.202616
28|46060|0.337|07.02.202617
29|46060|0.123|07.02.202618
30|46060|0|07.02.202619
31|46060|0|07.02.202620
32|46060|0|07.02.202621
33|46060|0|07.02.202622
34|46060|0|07.02.202623
35|46060|0|07.02.202624
36|46061|0|08.02.20261
37|46061|0|08.02.20262
38|46061|0|08.02.20263
39|46061|0|08.02.20264
40|46061|0|08.02.20265
41|46061|0|08.02.20266
42|46061|0|08.02.20267
43|46061|0|08.02.20268
44|46061|0|08.02.20269
45|46061|0.21|08.02.202610
46|46061|0.465|08.02.202611
47|46061|0.719|08.02.202612
48|46061|0.778|08.02.202613
49|46061|0.778|08.02.202614
50|46061|0.768|08.02.202615
51|46061|0.584|08.02.202616
52|46061|0.382|08.02.202617
53|46061|0.123|08.02.202618
54|46061|0|08.02.202619
55|46061|0|08.02.202620
56|46061|0|08.02.202621
57|46061|0|08.02.202622
58|46061|0|08.02.202623
59|46061|0|08.02.202624
60|46062|0|09.02.20261
61|46062|0|09.02.20262
62|46062|0|09.02.20263
63|46062|0|09.02.20264
64|46062|0|09.02.20265
65|46062|0|09.02.20266
66|46062|0|09.02.20267
67|46062|0|09.02.20268
68|46062|0|09.02.20269
69|46062|0.123|09.02.202610
70|46062|0.341|09.02.202611
71|46062|0.54|09.02.202612
72|46062|0.708|09.02.202613
73|46062|0.731|09.02.202614
74|46062|0.708|09.02.202615
75|46062|0.529|09.02.202616
76|46062|0.43|09.02.202617
77|46062|0.112|09.02.202618
78|46062|0|09.02.202619
79|46062|0|09.02.202620
80|46062|0|09.02.202621
81|46062|0|09.02.202622
82|46062|0|09.02.202623
83|46062|0|09.02.202624
84|46063|0|10.02.20261
85|46063|0|10.02.20262
86|46063|0|10.02.20263
87|46063|0|10.02.20264
88|46063|0|10.02.20265
89|46063|0|10.02.20266
90|46063|0|10.02.20267
91|46063|0|10.02.20268
92|46063|0|10.02.20269
93|46063|0.208|10.02.202610
94|46063|0.491|10.02.202611
95|46063|0.805|10.02.202612
96|46063|0.836|10.02.202613
97|46063|0.836|10.02.202614
98|46063|0.793|10.02.202615
99|46063|0.5580000000000001|10.02.202616
100|46063|0.333|10.02.202617
101|46063|0.112|10.02.202618
102|46063|0|10.02.202619
103|46063|0|10.02.202620
104|46063|0|10.02.202621
105|46063|0|10.02.202622
106|46063|0|10.02.202623
107|46063|0|10.02.202624
108|46064|0|11.02.20261
109|46064|0|11.02.20262
110|46064|0|11.02.20263
111|46064|0|11.02.20264
112|46064|0|11.02.20265
113|46064|0|11.02.20266
114|46064|0|11.02.20267
115|46064|0|11.02.20268
116|46064|0.015|11.02.20269
117|46064|0.316|11.02.202610
118|46064|0.74|11.02.202611
119|46064|1.011|11.02.202612
120|46064|1.217|11.02.202613
121|46064|1.206|11.02.202614
122|46064|1.169|11.02.202615
123|46064|1.12|11.02.202616
124|46064|0.757|11.02.202617
125|46064|0.352|11.02.202618
126|46064|0.023|11.02.202619
127|46064|0|11.02.202620
128|46064|0|11.02.202621
129|46064|0|11.02.202622
130|46064|0|11.02.202623
131|46064|0|11.02.202624
132|46065|0|12.02.20261
133|46065|0|12.02.20262
134|46065|0|12.02.20263
135|46065|0|12.02.20264
136|46065|0|12.02.20265
137|46065|0|12.02.20266
138|46065|0|12.02.20267
139|46065|0|12.02.20268
140|46065|0.015|12.02.20269
141|46065|0.326|12.02.202610
142|46065|0.802|12.02.202611
143|46065|1.005|12.02.202612
144|46065|1.212|12.02.202613
145|46065|1.166|12.02.202614
146|46065|1.176|12.02.202615
147|46065|1.005|12.02.202616
148|46065|0.766|12.02.202617
149|46065|0.315|12.02.202618
150|46065|0.019|12.02.202619
151|46065|0|12.02.202620
152|46065|0|12.02.202621
153|46065|0|12.02.202622
154|46065|0|12.02.202623
155|46065|0|12.02.202624
156|46066|0|13.02.20261
157|46066|0|13.02.20262
158|46066|0|13.02.20263
159|46066|0|13.02.20264
160|46066|0|13.02.20265
161|46066|0|13.02.20266
162|46066|0|13.02.20267
163|46066|0|13.02.20268
164|46066|0.014|13.02.20269
165|46066|0.254|13.02.202610
166|46066|0.582|13.02.202611
167|46066|0.774|13.02.202612
168|46066|0.97|13.02.202613
169|46066|0.909|13.02.202614
170|46066|0.748|13.02.202615
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $row = [int]$parts[0]
    $A = [double]$parts[1]
    $C = [double]$parts[2]
    $D = $parts[3]

    $ws.Cells.Item($row, 1).Value = $A
    $ws.Cells.Item($row, 3).Value = $C
    $ws.Cells.Item($row, 4).Value = $D
}

Write-Output "Updated $($lines.Count) rows"
